# Adapt column header formatting to respective input file names.
# - rename "<Column>_old" headers to "<Column>_FV2310"
# - rename "<Column>_new" headers to "<Column>_FV2404"
# - turn the sheet's data range into an Excel Table ("Table1")
# - freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$headers = @(
  "Segmentname_FV2310","Segmentgruppe_FV2310","Segment_FV2310","Datenelement_FV2310","Segment ID_FV2310",
  "Code_FV2310","Qualifier_FV2310","Beschreibung_FV2310","Bedingungsausdruck_FV2310","Bedingung_FV2310",
  "diff",
  "Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404",
  "Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn A1:U81 into an Excel Table (Table1), matching the existing header names.
$dataRange = $ws.Range("A1:U81")
$table = $ws.ListObjects.Add(1, $dataRange, [System.Type]::Missing, 1)
$table.Name = "Table1"
$table.TableStyle = ""

# Freeze the top (header) row.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
